$wb = $excel.ActiveWorkbook

$wsPIR = $wb.Worksheets.Item("PIR")
$wsPIR.Range("A86:A97").NumberFormat = "@"
$wsPIR.Range("E86:E97").NumberFormat = "@"
$wsPIR.Cells.Item(86, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(86, 2).Value = '17:09:07'
$wsPIR.Cells.Item(86, 3).Value = '17:00'
$wsPIR.Cells.Item(86, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(86, 5).Value = 'No Motion'
$wsPIR.Cells.Item(86, 6).Value = 'Inactive'
$wsPIR.Cells.Item(87, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(87, 2).Value = '17:09:08'
$wsPIR.Cells.Item(87, 3).Value = '17:00'
$wsPIR.Cells.Item(87, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(87, 5).Value = 'No Motion'
$wsPIR.Cells.Item(87, 6).Value = 'Inactive'
$wsPIR.Cells.Item(88, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(88, 2).Value = '17:09:09'
$wsPIR.Cells.Item(88, 3).Value = '17:00'
$wsPIR.Cells.Item(88, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(88, 5).Value = 'No Motion'
$wsPIR.Cells.Item(88, 6).Value = 'Inactive'
$wsPIR.Cells.Item(89, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(89, 2).Value = '17:09:14'
$wsPIR.Cells.Item(89, 3).Value = '17:00'
$wsPIR.Cells.Item(89, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(89, 5).Value = 'No Motion'
$wsPIR.Cells.Item(89, 6).Value = 'Inactive'
$wsPIR.Cells.Item(90, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(90, 2).Value = '17:09:19'
$wsPIR.Cells.Item(90, 3).Value = '17:00'
$wsPIR.Cells.Item(90, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(90, 5).Value = 'No Motion'
$wsPIR.Cells.Item(90, 6).Value = 'Inactive'
$wsPIR.Cells.Item(91, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(91, 2).Value = '17:09:24'
$wsPIR.Cells.Item(91, 3).Value = '17:00'
$wsPIR.Cells.Item(91, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(91, 5).Value = 'No Motion'
$wsPIR.Cells.Item(91, 6).Value = 'Inactive'
$wsPIR.Cells.Item(92, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(92, 2).Value = '17:09:29'
$wsPIR.Cells.Item(92, 3).Value = '17:00'
$wsPIR.Cells.Item(92, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(92, 5).Value = 'No Motion'
$wsPIR.Cells.Item(92, 6).Value = 'Inactive'
$wsPIR.Cells.Item(93, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(93, 2).Value = '17:09:35'
$wsPIR.Cells.Item(93, 3).Value = '17:00'
$wsPIR.Cells.Item(93, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(93, 5).Value = 'No Motion'
$wsPIR.Cells.Item(93, 6).Value = 'Inactive'
$wsPIR.Cells.Item(94, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(94, 2).Value = '17:09:39'
$wsPIR.Cells.Item(94, 3).Value = '17:00'
$wsPIR.Cells.Item(94, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(94, 5).Value = 'No Motion'
$wsPIR.Cells.Item(94, 6).Value = 'Inactive'
$wsPIR.Cells.Item(95, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(95, 2).Value = '17:09:44'
$wsPIR.Cells.Item(95, 3).Value = '17:00'
$wsPIR.Cells.Item(95, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(95, 5).Value = 'No Motion'
$wsPIR.Cells.Item(95, 6).Value = 'Inactive'
$wsPIR.Cells.Item(96, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(96, 2).Value = '17:09:50'
$wsPIR.Cells.Item(96, 3).Value = '17:00'
$wsPIR.Cells.Item(96, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(96, 5).Value = 'No Motion'
$wsPIR.Cells.Item(96, 6).Value = 'Inactive'
$wsPIR.Cells.Item(97, 1).Value = '2026-01-28'
$wsPIR.Cells.Item(97, 2).Value = '17:09:55'
$wsPIR.Cells.Item(97, 3).Value = '17:00'
$wsPIR.Cells.Item(97, 4).Value = 'Bathroom'
$wsPIR.Cells.Item(97, 5).Value = 'No Motion'
$wsPIR.Cells.Item(97, 6).Value = 'Inactive'

$wsHumidity = $wb.Worksheets.Item("Humidity")
$wsHumidity.Range("A86:A97").NumberFormat = "@"
$wsHumidity.Range("E86:E97").NumberFormat = "@"
$wsHumidity.Cells.Item(86, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(86, 2).Value = '17:09:07'
$wsHumidity.Cells.Item(86, 3).Value = '17:00'
$wsHumidity.Cells.Item(86, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(86, 5).Value = '87.5%'
$wsHumidity.Cells.Item(86, 6).Value = 'Active'
$wsHumidity.Cells.Item(87, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(87, 2).Value = '17:09:09'
$wsHumidity.Cells.Item(87, 3).Value = '17:00'
$wsHumidity.Cells.Item(87, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(87, 5).Value = '87.6%'
$wsHumidity.Cells.Item(87, 6).Value = 'Active'
$wsHumidity.Cells.Item(88, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(88, 2).Value = '17:09:13'
$wsHumidity.Cells.Item(88, 3).Value = '17:00'
$wsHumidity.Cells.Item(88, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(88, 5).Value = '87.6%'
$wsHumidity.Cells.Item(88, 6).Value = 'Active'
$wsHumidity.Cells.Item(89, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(89, 2).Value = '17:09:17'
$wsHumidity.Cells.Item(89, 3).Value = '17:00'
$wsHumidity.Cells.Item(89, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(89, 5).Value = '86.7%'
$wsHumidity.Cells.Item(89, 6).Value = 'Active'
$wsHumidity.Cells.Item(90, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(90, 2).Value = '17:09:21'
$wsHumidity.Cells.Item(90, 3).Value = '17:00'
$wsHumidity.Cells.Item(90, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(90, 5).Value = '87.6%'
$wsHumidity.Cells.Item(90, 6).Value = 'Active'
$wsHumidity.Cells.Item(91, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(91, 2).Value = '17:09:25'
$wsHumidity.Cells.Item(91, 3).Value = '17:00'
$wsHumidity.Cells.Item(91, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(91, 5).Value = '87.6%'
$wsHumidity.Cells.Item(91, 6).Value = 'Active'
$wsHumidity.Cells.Item(92, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(92, 2).Value = '17:09:29'
$wsHumidity.Cells.Item(92, 3).Value = '17:00'
$wsHumidity.Cells.Item(92, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(92, 5).Value = '86.7%'
$wsHumidity.Cells.Item(92, 6).Value = 'Active'
$wsHumidity.Cells.Item(93, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(93, 2).Value = '17:09:37'
$wsHumidity.Cells.Item(93, 3).Value = '17:00'
$wsHumidity.Cells.Item(93, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(93, 5).Value = '86.7%'
$wsHumidity.Cells.Item(93, 6).Value = 'Active'
$wsHumidity.Cells.Item(94, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(94, 2).Value = '17:09:45'
$wsHumidity.Cells.Item(94, 3).Value = '17:00'
$wsHumidity.Cells.Item(94, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(94, 5).Value = '87.6%'
$wsHumidity.Cells.Item(94, 6).Value = 'Active'
$wsHumidity.Cells.Item(95, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(95, 2).Value = '17:09:49'
$wsHumidity.Cells.Item(95, 3).Value = '17:00'
$wsHumidity.Cells.Item(95, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(95, 5).Value = '86.7%'
$wsHumidity.Cells.Item(95, 6).Value = 'Active'
$wsHumidity.Cells.Item(96, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(96, 2).Value = '17:09:53'
$wsHumidity.Cells.Item(96, 3).Value = '17:00'
$wsHumidity.Cells.Item(96, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(96, 5).Value = '87.7%'
$wsHumidity.Cells.Item(96, 6).Value = 'Active'
$wsHumidity.Cells.Item(97, 1).Value = '2026-01-28'
$wsHumidity.Cells.Item(97, 2).Value = '17:09:57'
$wsHumidity.Cells.Item(97, 3).Value = '17:00'
$wsHumidity.Cells.Item(97, 4).Value = 'Bathroom'
$wsHumidity.Cells.Item(97, 5).Value = '86.7%'
$wsHumidity.Cells.Item(97, 6).Value = 'Active'

$wsTemperature = $wb.Worksheets.Item("Temperature")
$wsTemperature.Range("A86:A97").NumberFormat = "@"
$wsTemperature.Range("E86:E97").NumberFormat = "@"
$wsTemperature.Cells.Item(86, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(86, 2).Value = '17:09:08'
$wsTemperature.Cells.Item(86, 3).Value = '17:00'
$wsTemperature.Cells.Item(86, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(86, 5).Value = '22.8C'
$wsTemperature.Cells.Item(86, 6).Value = 'Active'
$wsTemperature.Cells.Item(87, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(87, 2).Value = '17:09:09'
$wsTemperature.Cells.Item(87, 3).Value = '17:00'
$wsTemperature.Cells.Item(87, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(87, 5).Value = '22.8C'
$wsTemperature.Cells.Item(87, 6).Value = 'Active'
$wsTemperature.Cells.Item(88, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(88, 2).Value = '17:09:13'
$wsTemperature.Cells.Item(88, 3).Value = '17:00'
$wsTemperature.Cells.Item(88, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(88, 5).Value = '22.8C'
$wsTemperature.Cells.Item(88, 6).Value = 'Active'
$wsTemperature.Cells.Item(89, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(89, 2).Value = '17:09:17'
$wsTemperature.Cells.Item(89, 3).Value = '17:00'
$wsTemperature.Cells.Item(89, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(89, 5).Value = '22.9C'
$wsTemperature.Cells.Item(89, 6).Value = 'Active'
$wsTemperature.Cells.Item(90, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(90, 2).Value = '17:09:21'
$wsTemperature.Cells.Item(90, 3).Value = '17:00'
$wsTemperature.Cells.Item(90, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(90, 5).Value = '22.9C'
$wsTemperature.Cells.Item(90, 6).Value = 'Active'
$wsTemperature.Cells.Item(91, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(91, 2).Value = '17:09:25'
$wsTemperature.Cells.Item(91, 3).Value = '17:00'
$wsTemperature.Cells.Item(91, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(91, 5).Value = '22.8C'
$wsTemperature.Cells.Item(91, 6).Value = 'Active'
$wsTemperature.Cells.Item(92, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(92, 2).Value = '17:09:29'
$wsTemperature.Cells.Item(92, 3).Value = '17:00'
$wsTemperature.Cells.Item(92, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(92, 5).Value = '22.8C'
$wsTemperature.Cells.Item(92, 6).Value = 'Active'
$wsTemperature.Cells.Item(93, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(93, 2).Value = '17:09:37'
$wsTemperature.Cells.Item(93, 3).Value = '17:00'
$wsTemperature.Cells.Item(93, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(93, 5).Value = '22.8C'
$wsTemperature.Cells.Item(93, 6).Value = 'Active'
$wsTemperature.Cells.Item(94, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(94, 2).Value = '17:09:45'
$wsTemperature.Cells.Item(94, 3).Value = '17:00'
$wsTemperature.Cells.Item(94, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(94, 5).Value = '22.8C'
$wsTemperature.Cells.Item(94, 6).Value = 'Active'
$wsTemperature.Cells.Item(95, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(95, 2).Value = '17:09:49'
$wsTemperature.Cells.Item(95, 3).Value = '17:00'
$wsTemperature.Cells.Item(95, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(95, 5).Value = '22.8C'
$wsTemperature.Cells.Item(95, 6).Value = 'Active'
$wsTemperature.Cells.Item(96, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(96, 2).Value = '17:09:53'
$wsTemperature.Cells.Item(96, 3).Value = '17:00'
$wsTemperature.Cells.Item(96, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(96, 5).Value = '22.9C'
$wsTemperature.Cells.Item(96, 6).Value = 'Active'
$wsTemperature.Cells.Item(97, 1).Value = '2026-01-28'
$wsTemperature.Cells.Item(97, 2).Value = '17:09:57'
$wsTemperature.Cells.Item(97, 3).Value = '17:00'
$wsTemperature.Cells.Item(97, 4).Value = 'Bathroom'
$wsTemperature.Cells.Item(97, 5).Value = '22.8C'
$wsTemperature.Cells.Item(97, 6).Value = 'Active'

